$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (24) of kaspa buy data, run on 2026-01-23.
# The date column is stored as plain text (matching the existing
# convention used by the other date rows), so force a text format before
# assigning the value to stop Excel from auto-converting the "MM/DD/YYYY"
# string into a date serial number, then restore the default "Normal"
# style so the cell doesn't end up with a lingering custom number format.
$dateCell = $ws.Range("A24")
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/23/2026"
$dateCell.Style = "Normal"

$ws.Range("B24").Value = 1133.198
$ws.Range("C24").Value = 0.04368168669552892
$ws.Range("D24").Value = 50
